$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2 through 10 (9 data rows), columns A-T.
# A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster
# then E..T numeric metrics.
# Values below are listed per data row (row2..row10) in column order A,B,C,D,E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T

$rows = @(
    @{A="ECs";  B="Wnt5a"; C="Ror2"; D="ECs";  E=1; F=0.3333333333333333; G=0.04240166666666667; H=0.127205;          I=0.004737824075193167; J=0.004737824075193167; K=1; L=0.3333333333333333; M=0.1389836666666667; N=0.416951;  O=0.01758007579884798; P=0.01758007579884798; Q=0.005893139106111112; R=0.053038251955;      S=0.00008329130636350272; T=0.00008329130636350272},
    @{A="ECs";  B="Wnt5a"; C="Ror2"; D="FAPs"; E=1; F=0.3333333333333333; G=0.04240166666666667; H=0.127205;          I=0.004737824075193167; J=0.004737824075193167; K=3; L=1;                  M=7.377589;           N=22.132767; O=0.9331929207466617;  P=0.9331929207466617;  Q=0.3128220695816667;   R=2.815398626235;      S=0.004421303886713363;  T=0.004421303886713363},
    @{A="ECs";  B="Wnt5a"; C="Ror2"; D="sCs";  E=1; F=0.3333333333333333; G=0.04240166666666667; H=0.127205;          I=0.004737824075193167; J=0.004737824075193167; K=2; L=0.6666666666666666; M=0.3891763333333333; N=1.167529;  O=0.04922700345449031; P=0.04922700345449031; Q=0.01650172516055556;  R=0.148515526445;      S=0.0002332288821163014; T=0.0002332288821163014},
    @{A="FAPs"; B="Wnt5a"; C="Ror2"; D="ECs";  E=3; F=1;                  G=8.775005999999999;  H=26.325018;          I=0.9804905786745289;  J=0.9804905786745289;  K=1; L=0.3333333333333333; M=0.1389836666666667; N=0.416951;  O=0.01758007579884798; P=0.01758007579884798; Q=1.219582508902;       R=10.976242580118;     S=0.01723709869315454;  T=0.01723709869315454},
    @{A="FAPs"; B="Wnt5a"; C="Ror2"; D="FAPs"; E=3; F=1;                  G=8.775005999999999;  H=26.325018;          I=0.9804905786745289;  J=0.9804905786745289;  K=3; L=1;                  M=7.377589;           N=22.132767; O=0.9331929207466617;  P=0.9331929207466617;  Q=64.73838774053399;    R=582.6454896648061;   S=0.9149868668778681;    T=0.9149868668778681},
    @{A="FAPs"; B="Wnt5a"; C="Ror2"; D="sCs";  E=3; F=1;                  G=8.775005999999999;  H=26.325018;          I=0.9804905786745289;  J=0.9804905786745289;  K=2; L=0.6666666666666666; M=0.3891763333333333; N=1.167529;  O=0.04922700345449031; P=0.04922700345449031; Q=3.415024660058;       R=30.735221940522;     S=0.04826661310350624;   T=0.04826661310350624},
    @{A="sCs";  B="Wnt5a"; C="Ror2"; D="ECs";  E=1; F=0.3333333333333333; G=0.1322;             H=0.3966;             I=0.01477159725027798;  J=0.01477159725027797; K=1; L=0.3333333333333333; M=0.1389836666666667; N=0.416951;  O=0.01758007579884798; P=0.01758007579884798; Q=0.01837364073333334;  R=0.1653627666;        S=0.0002596857993299413; T=0.0002596857993299412},
    @{A="sCs";  B="Wnt5a"; C="Ror2"; D="FAPs"; E=1; F=0.3333333333333333; G=0.1322;             H=0.3966;             I=0.01477159725027798;  J=0.01477159725027797; K=3; L=1;                  M=7.377589;           N=22.132767; O=0.9331929207466617;  P=0.9331929207466617;  Q=0.9753172658000001;   R=8.777855392200001;   S=0.01378474998208026;   T=0.01378474998208026},
    @{A="sCs";  B="Wnt5a"; C="Ror2"; D="sCs";  E=1; F=0.3333333333333333; G=0.1322;             H=0.3966;             I=0.01477159725027798;  J=0.01477159725027797; K=2; L=0.6666666666666666; M=0.3891763333333333; N=1.167529;  O=0.04922700345449031; P=0.04922700345449031; Q=0.05144911126666667;  R=0.4630420014;        S=0.0007271614688677735; T=0.0007271614688677733}
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Write column-by-column (all rows for column A, then all rows for column B, etc.)
# so that new shared strings are registered in the same order Excel would
# encounter them (by column, top-to-bottom).
foreach ($col in $columns) {
    $rowIndex = 2
    foreach ($row in $rows) {
        $cellRef = "$col$rowIndex"
        $ws.Range($cellRef).Value = $row[$col]
        $rowIndex++
    }
}
